$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell F1 ("Remark") backed by the shared string that used to hold "Tom"
$ws.Range("F1").Value = "Remark"

# A2 used to be the literal shared string "Tom"; it now derives its value from a
# formula that concatenates the new "Remark" header with the existing F2 value.
$ws.Range("A2").Formula = "=CONCATENATE(F1,""_"",F2)"

# Move the active selection from C3 to A2 to match the saved view state.
$ws.Range("A2").Select() | Out-Null
